$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from "Property1" to "DataNode"
$ws.Name = "DataNode"

# Move the active selection from A9 to C41 (within the frozen "bottomLeft" pane)
$ws.Range("C41").Select()
